$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# NOTE: this PowerShell subset does not bind named (-Param value) arguments
# inside a param() block, so every helper below uses plain positional
# parameters only.
function AddRow($NewRow, $SrcRow, $Timestamp, $Name, $Roll, $Profile, $Branch, $Score, $HlTarget, $HlDisplay) {
    # Duplicate the formatting (styles/fonts/borders/wrap/number-format) from
    # the template row so the new row reuses the workbook's existing style
    # indices instead of minting new ones.
    $ws.Range("A$SrcRow`:F$SrcRow").Copy()
    $ws.Range("A$NewRow`:F$NewRow").PasteSpecial($xlPasteFormats)

    $ws.Cells.Item($NewRow, 1).Value = $Timestamp
    $ws.Cells.Item($NewRow, 2).Value = $Name
    $ws.Cells.Item($NewRow, 3).Value = $Roll
    $ws.Cells.Item($NewRow, 4).Value = $Profile
    $ws.Cells.Item($NewRow, 5).Value = $Branch
    $ws.Cells.Item($NewRow, 6).Value = $Score

    if ($HlTarget -ne "") {
        if ($HlDisplay -ne "") {
            $ws.Hyperlinks.Add($ws.Cells.Item($NewRow, 4), $HlTarget, "", "", $HlDisplay)
            $ws.Cells.Item($NewRow, 4).Value = $Profile
        } else {
            $ws.Hyperlinks.Add($ws.Cells.Item($NewRow, 4), $HlTarget)
        }
        # Hyperlinks.Add stamps the built-in "Hyperlink" cell style; restore
        # the sheet's own hyperlinked-cell format (border + wrap) to match
        # the template row's D cell.
        $ws.Range("D$SrcRow").Copy()
        $ws.Range("D$NewRow").PasteSpecial($xlPasteFormats)
    }

    $ws.Rows.Item($NewRow).RowHeight = $ws.Rows.Item($SrcRow).RowHeight
}

# Rows 67-73 duplicate existing rows 60-66 verbatim (same timestamps, names,
# roll numbers, profile links, and branches).
AddRow 67 60 45384.440335648149 "Manya Gupta" "B23154" "https://www.beecrowd.com.br/judge/en/profile/949111" "CSE" 0 "https://www.beecrowd.com.br/judge/en/profile/949111" ""

AddRow 68 61 45384.446770833332 "Pranab Ray" "B23169" "https://www.beecrowd.com.br/judge/en/profile/942674" "CSE" 0 "" ""

AddRow 69 62 45384.469606481478 "Shubhankit Singh" "B23387" "https://www.beecrowd.com.br/judge/en/profile/948383" "MSE" 0 "" ""

AddRow 70 63 45384.486840277779 "Aaryan Tyagi" "B23420" "https://www.beecrowd.com.br/judge/en/profile/949129" "ME" 0 "https://www.beecrowd.com.br/judge/en/profile/949129" ""

AddRow 71 64 45384.522835648146 "Ansh Attre" "b23191" "https://www.beecrowd.com.br/judge/en/profile/948707" "DSE" 0 "" ""

AddRow 72 65 45384.52925925926 "Arpita Kumari" "B23249" "https://www.beecrowd.com.br/judge/en/profile/949137" "EE" 0 "" ""

AddRow 73 66 45384.53224537037 "Arka" "B23120" "https://www.beecrowd.com.br/judge/en/profile/948169" "CSE" 0 "" ""

# Row 74 is a genuinely new submission (Rohit), formatted like the other
# hyperlinked rows, with a cached display string that differs from the
# underlying address (mirrors the source file's stored hyperlink text).
AddRow 74 60 45384.410578703704 "ROHIT" "B23175" "https://www.beecrowd.com.br/judge/en/profile/948713" "CSE" 0 "https://www.beecrowd.com.br/judge/en/profile/948713" "https://www.beecrowd.com.br/judge/en/users/basic-info"

# Keep the view/selection state consistent with the appended rows, matching
# Excel's own behaviour when the used range grows past the visible window.
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$ws.Range("F78").Select()

Write-Output "Appended rows 67-74"
